# "improvements to tech review"
#
# 1) Slide 2 ("Sally the Scientist, Activist"): add a bullet "OS agnostic"
#    right after "Novice with web service technology".
# 2) Slide 3 ("Web service technologies"): add a level-1 bullet
#    "Available on all major OSes" (as two runs: "Available on " /
#    "all major OSes") right after "Minimal install requirements".

$p = $ppt.ActivePresentation

# --- Slide 2: content placeholder ---
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$paras2 = $tr2.Paragraphs()
$idx2 = -1
for ($i = 1; $i -le $paras2.Count; $i++) {
    if ($tr2.Paragraphs($i, 1).Text -like "*Novice with web service technology*") {
        $idx2 = $i
    }
}
if ($idx2 -le 0) { $idx2 = $paras2.Count }
$target2 = $tr2.Paragraphs($idx2, 1)
[void]$target2.InsertAfter("`rOS agnostic")

# --- Slide 3: content placeholder ---
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange
$paras3 = $tr3.Paragraphs()
$idx3 = -1
for ($i = 1; $i -le $paras3.Count; $i++) {
    if ($tr3.Paragraphs($i, 1).Text -like "*Minimal install requirements*") {
        $idx3 = $i
    }
}
if ($idx3 -le 0) { $idx3 = 4 }
$target3 = $tr3.Paragraphs($idx3, 1)
[void]$target3.InsertAfter("`rAvailable on all major OSes")

$newPara3 = $tr3.Paragraphs($idx3 + 1, 1)
# Split the new paragraph into two runs: "Available on " + "all major OSes"
$secondRun3 = $newPara3.Characters(14, 14)
$secondRun3.Text = $secondRun3.Text
